# Vendors p0.xlsx - "successfully writes expense sheet and updates stock"
#
# 1) Update remaining Inventory stock counts (post-shopping-list fulfillment).
# 2) Populate the (previously empty) Expenses sheet with a Client/Total
#    summary computed from Shopping List quantities x Inventory price/discount.

$wb = $excel.ActiveWorkbook

# ---- Inventory: refresh Stock column (D) ----
$ws1 = $wb.Worksheets.Item("Inventory")
$ws1.Range("D2").Value  = 4     # Socks
$ws1.Range("D3").Value  = 76    # Bananas
$ws1.Range("D5").Value  = 30    # Oranges
$ws1.Range("D6").Value  = 66    # Nesquik
$ws1.Range("D7").Value  = 80    # Candy
$ws1.Range("D8").Value  = 2     # Guitar
$ws1.Range("D10").Value = 2     # T.V

# ---- Expenses: write Client / Total report ----
$ws3 = $wb.Worksheets.Item("Expenses")
$ws3.Range("A1").Value = "Client"
$ws3.Range("B1").Value = "Total"
$ws3.Range("A2").Value = "Bob"
$ws3.Range("B2").Value = 211.25
$ws3.Range("A3").Value = "Webb"
$ws3.Range("B3").Value = 12.1
$ws3.Range("A4").Value = "Joe"
$ws3.Range("B4").Value = 78.75

# ---- Restore the cursor's final resting position on Inventory ----
[void]$ws1.Range("D28").Select()
